$d = $word.ActiveDocument

# --- Insert three new bullet paragraphs before the
# "Developed and deployed custom analytical tools..." bullet under the
# Siege Analytics / "Advanced Data Analysis and Statistical Modeling" section.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Developed and deployed custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering*") {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $insertionPoint = $targetPara.Range.Duplicate
    $insertionPoint.Collapse(1)
    $newBullets = "• Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data`r" + `
                  "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters`r" + `
                  "• Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts`r"
    $insertionPoint.InsertBefore($newBullets)
}

# --- Remove the now-superseded bullet about campaign finance fraud detection.
$removePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets*") {
        $removePara = $p
        break
    }
}

if ($removePara -ne $null) {
    $removePara.Range.Delete()
}
